# Scheduled-runner update: refresh Leve profit/price figures across sheets.
# Applies updated currentAveragePrice* / LevePrice* / LeveProfit* figures
# (columns H-N) for a set of rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# sheets, as produced by the market-data refresh job.

$wb = $excel.ActiveWorkbook

function Set-RowValues {
    param($Sheet, $Row, $Values)
    foreach ($col in $Values.Keys) {
        $addr = "$col$Row"
        $Sheet.Range($addr).Value = $Values[$col]
    }
}

# ---------------------------------------------------------------------------
# ALC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

Set-RowValues $ws 12 @{
    "H" = 729.05884
    "I" = 634.4
    "K" = 634.4
    "M" = -464.4
}

Set-RowValues $ws 70 @{
    "H" = 950
    "I" = 600
    "J" = 1300
    "K" = 1800
    "L" = 3900
    "M" = -1530
    "N" = -4440
}

Set-RowValues $ws 73 @{
    "H" = 950
    "I" = 600
    "J" = 1300
    "K" = 1800
    "L" = 3900
    "M" = -864
    "N" = -5772
}

Set-RowValues $ws 112 @{
    "H" = 1892.6727
    "J" = 1916.6111
    "L" = 5749.8333
    "N" = -7965.8333
}

Set-RowValues $ws 137 @{
    "H" = 85459
    "I" = 2016.5
    "K" = 6049.5
    "M" = -3499.5
}

Set-RowValues $ws 138 @{
    "H" = 3037.4268
    "I" = 4120.778
    "J" = 2732.7344
    "K" = 12362.334
    "L" = 8198.2032
    "M" = -7222.334000000001
    "N" = -18478.2032
}

# ---------------------------------------------------------------------------
# ARM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

Set-RowValues $ws 74 @{
    "H" = 876.55316
    "I" = 707.0476
    "K" = 707.0476
    "M" = 166.9524
}

Set-RowValues $ws 77 @{
    "H" = 876.55316
    "I" = 707.0476
    "K" = 3535.238
    "M" = 832.7620000000002
}

Set-RowValues $ws 122 @{
    "H" = 1728.8
    "I" = 1657.6364
    "J" = 1924.5
    "K" = 4972.9092
    "L" = 5773.5
    "M" = -2522.9092
    "N" = -10673.5
}

Set-RowValues $ws 132 @{
    "H" = 1860.2273
    "I" = 1482.8605
    "K" = 4448.5815
    "M" = -1918.5815
}

# ---------------------------------------------------------------------------
# BSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

Set-RowValues $ws 134 @{
    "H" = 4554.689
    "I" = 4660.0835
    "K" = 13980.2505
    "M" = -11445.2505
}

# ---------------------------------------------------------------------------
# CRP
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

Set-RowValues $ws 4 @{
    "H" = 4476
    "J" = 5370
    "L" = 5370
    "N" = -5594
}

Set-RowValues $ws 62 @{
    "H" = 3533.8333
    "I" = 2765.6667
    "J" = 4302
    "K" = 2765.6667
    "L" = 4302
    "M" = -2141.6667
    "N" = -5550
}

Set-RowValues $ws 65 @{
    "H" = 3533.8333
    "I" = 2765.6667
    "J" = 4302
    "K" = 13828.3335
    "L" = 21510
    "M" = -10708.3335
    "N" = -27750
}

Set-RowValues $ws 70 @{
    "H" = 50030
    "J" = 50030
    "L" = 50030
    "N" = -50660
}

Set-RowValues $ws 73 @{
    "H" = 50030
    "J" = 50030
    "L" = 50030
    "N" = -52214
}

Set-RowValues $ws 105 @{
    "H" = 611.6111
    "I" = 539
    "K" = 539
    "M" = 1208
}

Set-RowValues $ws 132 @{
    "H" = 1499.8462
    "I" = 1135.7222
    "K" = 3407.1666
    "M" = -877.1665999999996
}

Set-RowValues $ws 134 @{
    "H" = 1419.8148
    "I" = 1278.8636
    "K" = 3836.5908
    "M" = -1301.5908
}

# ---------------------------------------------------------------------------
# CUL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

Set-RowValues $ws 4 @{
    "H" = 1022614.6
    "I" = 1166865.8
    "K" = 3500597.4
    "M" = -3500485.4
}

Set-RowValues $ws 122 @{
    "H" = 1574.5
    "J" = 1728
    "L" = 15552
    "N" = -20452
}

# ---------------------------------------------------------------------------
# GSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

Set-RowValues $ws 5 @{
    "H" = 3833.3333
    "I" = 1000
    "K" = 1000
    "M" = -888
}

Set-RowValues $ws 97 @{
    "H" = 2366.4666
    "I" = 2454.3635
    "K" = 2454.3635
    "M" = -1958.3635
}

Set-RowValues $ws 102 @{
    "H" = 2909.1428
    "J" = 1184
    "L" = 1184
    "N" = -4428
}

Set-RowValues $ws 122 @{
    "H" = 1266
    "I" = 761.6
    "J" = 2527
    "K" = 2284.8
    "L" = 7581
    "M" = 165.1999999999998
    "N" = -12481
}

Set-RowValues $ws 132 @{
    "H" = 942191.7
    "I" = 1286542.9
    "J" = 3052
    "K" = 3859628.7
    "L" = 9156
    "M" = -3857098.7
    "N" = -14216
}

# ---------------------------------------------------------------------------
# LTW
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

Set-RowValues $ws 2 @{
    "H" = 354558.8
    "I" = 461576.94
    "J" = 6750
    "K" = 461576.94
    "L" = 6750
    "M" = -461464.94
    "N" = -6974
}

Set-RowValues $ws 40 @{
    "H" = 19790.6
    "J" = 10983.667
    "L" = 10983.667
    "N" = -11255.667
}

Set-RowValues $ws 122 @{
    "H" = 2998.5715
    "I" = 2798
    "J" = 3500
    "K" = 8394
    "L" = 10500
    "M" = -5944
    "N" = -15400
}

Set-RowValues $ws 130 @{
    "H" = 13214
    "J" = 13214
    "L" = 13214
    "N" = -23254
}

# ---------------------------------------------------------------------------
# WVR
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

Set-RowValues $ws 70 @{
    "H" = 47103.5
    "J" = 47103.5
    "L" = 47103.5
    "N" = -47733.5
}

Set-RowValues $ws 73 @{
    "H" = 47103.5
    "J" = 47103.5
    "L" = 47103.5
    "N" = -49287.5
}

Set-RowValues $ws 132 @{
    "H" = 1785.5682
    "I" = 1615.5
    "J" = 2862.6667
    "K" = 4846.5
    "L" = 8588.000100000001
    "M" = -2316.5
    "N" = -13648.0001
}
